# The workbook contains 14 worksheets (tabs "18" down to "5"), each with a
# single large text cell in B2 holding the verbatim console output of a
# statsmodels OLS regression summary. The summary text embeds the
# date/time the regression was run:
#
#     Date:                Sat, 28 Dec 2019   Prob (F-statistic): ...
#     Time:                        20:59:40   Log-Likelihood:     ...
#
# The model was simply re-run the next day, so every sheet's Date/Time
# stamp moves from Sat 28 Dec 2019 20:59:40 to Sun 29 Dec 2019 16:11:07 -
# none of the actual statistics changed. Walk every worksheet and patch
# cell B2 in place.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -ne $null -and $text -like "*Sat, 28 Dec 2019*") {
        $text = $text -replace "Sat, 28 Dec 2019", "Sun, 29 Dec 2019"
        $text = $text -replace "20:59:40", "16:11:07"
        $cell.Value = $text
    }
}

Write-Output "done"
